# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 135
$wsExhibit.Range("F6").Value = 9157
$wsExhibit.Range("F7").Value = 831
$wsExhibit.Range("F10").Value = 1057
$wsExhibit.Range("F15").Value = 353
$wsExhibit.Range("F18").Value = 1187

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 135
$wsAll.Range("F8").Value = 9157
$wsAll.Range("F9").Value = 831
$wsAll.Range("F12").Value = 1057
$wsAll.Range("F17").Value = 353
$wsAll.Range("F20").Value = 1187
